# Adapt tests to control version
# Add a "version" column to the "settings" sheet with value 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

$ws.Range("C1").Value = "version"
$ws.Range("C2").Value = 1
